# Auto-generated edit script: updates Distance/Size-coded condition
# labels and filenames in the shared strings used throughout the sheet
# (D80->D86, D64->D69, D51->D55, S30->S31), matching the regenerated order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "Fixation_D69_l.png"
$ws.Range("E2").Value2 = "Fixation_D69_r.png"
$ws.Range("H2").Value2 = "D69"
$ws.Range("D3").Value2 = "Fixation_D86_l.png"
$ws.Range("E3").Value2 = "Fixation_D86_r.png"
$ws.Range("H3").Value2 = "D86"
$ws.Range("B4").Value2 = "Face17_D86_S20"
$ws.Range("D4").Value2 = "Face17_D86_S20_l.png"
$ws.Range("E4").Value2 = "Face17_D86_S20_r.png"
$ws.Range("H4").Value2 = "D86"
$ws.Range("D5").Value2 = "Fixation_D86_l.png"
$ws.Range("E5").Value2 = "Fixation_D86_r.png"
$ws.Range("H5").Value2 = "D86"
$ws.Range("B6").Value2 = "Face07_D86_S31"
$ws.Range("D6").Value2 = "Face07_D86_S31_l.png"
$ws.Range("E6").Value2 = "Face07_D86_S31_r.png"
$ws.Range("H6").Value2 = "D86"
$ws.Range("J6").Value2 = "S31"
$ws.Range("D7").Value2 = "Fixation_D86_l.png"
$ws.Range("E7").Value2 = "Fixation_D86_r.png"
$ws.Range("H7").Value2 = "D86"
$ws.Range("B8").Value2 = "Face06_D86_S25"
$ws.Range("D8").Value2 = "Face06_D86_S25_l.png"
$ws.Range("E8").Value2 = "Face06_D86_S25_r.png"
$ws.Range("H8").Value2 = "D86"
$ws.Range("D9").Value2 = "Fixation_D86_l.png"
$ws.Range("E9").Value2 = "Fixation_D86_r.png"
$ws.Range("H9").Value2 = "D86"
$ws.Range("B10").Value2 = "Face09_D86_S20"
$ws.Range("D10").Value2 = "Face09_D86_S20_l.png"
$ws.Range("E10").Value2 = "Face09_D86_S20_r.png"
$ws.Range("H10").Value2 = "D86"
$ws.Range("D11").Value2 = "Fixation_D86_l.png"
$ws.Range("E11").Value2 = "Fixation_D86_r.png"
$ws.Range("H11").Value2 = "D86"
$ws.Range("B12").Value2 = "Face10_D86_S25"
$ws.Range("D12").Value2 = "Face10_D86_S25_l.png"
$ws.Range("E12").Value2 = "Face10_D86_S25_r.png"
$ws.Range("H12").Value2 = "D86"
$ws.Range("D13").Value2 = "Fixation_D86_l.png"
$ws.Range("E13").Value2 = "Fixation_D86_r.png"
$ws.Range("H13").Value2 = "D86"
$ws.Range("B14").Value2 = "Face01_D86_S31"
$ws.Range("D14").Value2 = "Face01_D86_S31_l.png"
$ws.Range("E14").Value2 = "Face01_D86_S31_r.png"
$ws.Range("H14").Value2 = "D86"
$ws.Range("J14").Value2 = "S31"
$ws.Range("D15").Value2 = "Fixation_D86_l.png"
$ws.Range("E15").Value2 = "Fixation_D86_r.png"
$ws.Range("H15").Value2 = "D86"
$ws.Range("D16").Value2 = "Fixation_D69_l.png"
$ws.Range("E16").Value2 = "Fixation_D69_r.png"
$ws.Range("H16").Value2 = "D69"
$ws.Range("B17").Value2 = "Face14_D69_S25"
$ws.Range("D17").Value2 = "Face14_D69_S25_l.png"
$ws.Range("E17").Value2 = "Face14_D69_S25_r.png"
$ws.Range("H17").Value2 = "D69"
$ws.Range("D18").Value2 = "Fixation_D69_l.png"
$ws.Range("E18").Value2 = "Fixation_D69_r.png"
$ws.Range("H18").Value2 = "D69"
$ws.Range("B19").Value2 = "Face02_D69_S25"
$ws.Range("D19").Value2 = "Face02_D69_S25_l.png"
$ws.Range("E19").Value2 = "Face02_D69_S25_r.png"
$ws.Range("H19").Value2 = "D69"
$ws.Range("D20").Value2 = "Fixation_D69_l.png"
$ws.Range("E20").Value2 = "Fixation_D69_r.png"
$ws.Range("H20").Value2 = "D69"
$ws.Range("B21").Value2 = "Face05_D69_S31"
$ws.Range("D21").Value2 = "Face05_D69_S31_l.png"
$ws.Range("E21").Value2 = "Face05_D69_S31_r.png"
$ws.Range("H21").Value2 = "D69"
$ws.Range("J21").Value2 = "S31"
$ws.Range("D22").Value2 = "Fixation_D69_l.png"
$ws.Range("E22").Value2 = "Fixation_D69_r.png"
$ws.Range("H22").Value2 = "D69"
$ws.Range("B23").Value2 = "Face03_D69_S31"
$ws.Range("D23").Value2 = "Face03_D69_S31_l.png"
$ws.Range("E23").Value2 = "Face03_D69_S31_r.png"
$ws.Range("H23").Value2 = "D69"
$ws.Range("J23").Value2 = "S31"
$ws.Range("D24").Value2 = "Fixation_D69_l.png"
$ws.Range("E24").Value2 = "Fixation_D69_r.png"
$ws.Range("H24").Value2 = "D69"
$ws.Range("B25").Value2 = "Face16_D69_S20"
$ws.Range("D25").Value2 = "Face16_D69_S20_l.png"
$ws.Range("E25").Value2 = "Face16_D69_S20_r.png"
$ws.Range("H25").Value2 = "D69"
$ws.Range("D26").Value2 = "Fixation_D69_l.png"
$ws.Range("E26").Value2 = "Fixation_D69_r.png"
$ws.Range("H26").Value2 = "D69"
$ws.Range("B27").Value2 = "Face15_D69_S20"
$ws.Range("D27").Value2 = "Face15_D69_S20_l.png"
$ws.Range("E27").Value2 = "Face15_D69_S20_r.png"
$ws.Range("H27").Value2 = "D69"
$ws.Range("D28").Value2 = "Fixation_D69_l.png"
$ws.Range("E28").Value2 = "Fixation_D69_r.png"
$ws.Range("H28").Value2 = "D69"
$ws.Range("D29").Value2 = "Fixation_D55_l.png"
$ws.Range("E29").Value2 = "Fixation_D55_r.png"
$ws.Range("H29").Value2 = "D55"
$ws.Range("B30").Value2 = "Face12_D55_S20"
$ws.Range("D30").Value2 = "Face12_D55_S20_l.png"
$ws.Range("E30").Value2 = "Face12_D55_S20_r.png"
$ws.Range("H30").Value2 = "D55"
$ws.Range("D31").Value2 = "Fixation_D55_l.png"
$ws.Range("E31").Value2 = "Fixation_D55_r.png"
$ws.Range("H31").Value2 = "D55"
$ws.Range("B32").Value2 = "Face04_D55_S20"
$ws.Range("D32").Value2 = "Face04_D55_S20_l.png"
$ws.Range("E32").Value2 = "Face04_D55_S20_r.png"
$ws.Range("H32").Value2 = "D55"
$ws.Range("D33").Value2 = "Fixation_D55_l.png"
$ws.Range("E33").Value2 = "Fixation_D55_r.png"
$ws.Range("H33").Value2 = "D55"
$ws.Range("B34").Value2 = "Face18_D55_S25"
$ws.Range("D34").Value2 = "Face18_D55_S25_l.png"
$ws.Range("E34").Value2 = "Face18_D55_S25_r.png"
$ws.Range("H34").Value2 = "D55"
$ws.Range("D35").Value2 = "Fixation_D55_l.png"
$ws.Range("E35").Value2 = "Fixation_D55_r.png"
$ws.Range("H35").Value2 = "D55"
$ws.Range("B36").Value2 = "Face08_D55_S25"
$ws.Range("D36").Value2 = "Face08_D55_S25_l.png"
$ws.Range("E36").Value2 = "Face08_D55_S25_r.png"
$ws.Range("H36").Value2 = "D55"
$ws.Range("D37").Value2 = "Fixation_D55_l.png"
$ws.Range("E37").Value2 = "Fixation_D55_r.png"
$ws.Range("H37").Value2 = "D55"
$ws.Range("B38").Value2 = "Face11_D55_S31"
$ws.Range("D38").Value2 = "Face11_D55_S31_l.png"
$ws.Range("E38").Value2 = "Face11_D55_S31_r.png"
$ws.Range("H38").Value2 = "D55"
$ws.Range("J38").Value2 = "S31"
$ws.Range("D39").Value2 = "Fixation_D55_l.png"
$ws.Range("E39").Value2 = "Fixation_D55_r.png"
$ws.Range("H39").Value2 = "D55"
$ws.Range("B40").Value2 = "Face13_D55_S31"
$ws.Range("D40").Value2 = "Face13_D55_S31_l.png"
$ws.Range("E40").Value2 = "Face13_D55_S31_r.png"
$ws.Range("H40").Value2 = "D55"
$ws.Range("J40").Value2 = "S31"
$ws.Range("D41").Value2 = "Fixation_D55_l.png"
$ws.Range("E41").Value2 = "Fixation_D55_r.png"
$ws.Range("H41").Value2 = "D55"
$ws.Range("D42").Value2 = "Fixation_D69_l.png"
$ws.Range("E42").Value2 = "Fixation_D69_r.png"
$ws.Range("H42").Value2 = "D69"
$ws.Range("B43").Value2 = "Face07_D69_S25"
$ws.Range("D43").Value2 = "Face07_D69_S25_l.png"
$ws.Range("E43").Value2 = "Face07_D69_S25_r.png"
$ws.Range("H43").Value2 = "D69"
$ws.Range("D44").Value2 = "Fixation_D69_l.png"
$ws.Range("E44").Value2 = "Fixation_D69_r.png"
$ws.Range("H44").Value2 = "D69"
$ws.Range("B45").Value2 = "Face04_D69_S31"
$ws.Range("D45").Value2 = "Face04_D69_S31_l.png"
$ws.Range("E45").Value2 = "Face04_D69_S31_r.png"
$ws.Range("H45").Value2 = "D69"
$ws.Range("J45").Value2 = "S31"
$ws.Range("D46").Value2 = "Fixation_D69_l.png"
$ws.Range("E46").Value2 = "Fixation_D69_r.png"
$ws.Range("H46").Value2 = "D69"
$ws.Range("B47").Value2 = "Face11_D69_S31"
$ws.Range("D47").Value2 = "Face11_D69_S31_l.png"
$ws.Range("E47").Value2 = "Face11_D69_S31_r.png"
$ws.Range("H47").Value2 = "D69"
$ws.Range("J47").Value2 = "S31"
$ws.Range("D48").Value2 = "Fixation_D69_l.png"
$ws.Range("E48").Value2 = "Fixation_D69_r.png"
$ws.Range("H48").Value2 = "D69"
$ws.Range("B49").Value2 = "Face03_D69_S25"
$ws.Range("D49").Value2 = "Face03_D69_S25_l.png"
$ws.Range("E49").Value2 = "Face03_D69_S25_r.png"
$ws.Range("H49").Value2 = "D69"
$ws.Range("D50").Value2 = "Fixation_D69_l.png"
$ws.Range("E50").Value2 = "Fixation_D69_r.png"
$ws.Range("H50").Value2 = "D69"
$ws.Range("B51").Value2 = "Face17_D69_S20"
$ws.Range("D51").Value2 = "Face17_D69_S20_l.png"
$ws.Range("E51").Value2 = "Face17_D69_S20_r.png"
$ws.Range("H51").Value2 = "D69"
$ws.Range("D52").Value2 = "Fixation_D69_l.png"
$ws.Range("E52").Value2 = "Fixation_D69_r.png"
$ws.Range("H52").Value2 = "D69"
$ws.Range("B53").Value2 = "Face02_D69_S20"
$ws.Range("D53").Value2 = "Face02_D69_S20_l.png"
$ws.Range("E53").Value2 = "Face02_D69_S20_r.png"
$ws.Range("H53").Value2 = "D69"
$ws.Range("D54").Value2 = "Fixation_D69_l.png"
$ws.Range("E54").Value2 = "Fixation_D69_r.png"
$ws.Range("H54").Value2 = "D69"
$ws.Range("D55").Value2 = "Fixation_D86_l.png"
$ws.Range("E55").Value2 = "Fixation_D86_r.png"
$ws.Range("H55").Value2 = "D86"
$ws.Range("B56").Value2 = "Face05_D86_S20"
$ws.Range("D56").Value2 = "Face05_D86_S20_l.png"
$ws.Range("E56").Value2 = "Face05_D86_S20_r.png"
$ws.Range("H56").Value2 = "D86"
$ws.Range("D57").Value2 = "Fixation_D86_l.png"
$ws.Range("E57").Value2 = "Fixation_D86_r.png"
$ws.Range("H57").Value2 = "D86"
$ws.Range("B58").Value2 = "Face13_D86_S20"
$ws.Range("D58").Value2 = "Face13_D86_S20_l.png"
$ws.Range("E58").Value2 = "Face13_D86_S20_r.png"
$ws.Range("H58").Value2 = "D86"
$ws.Range("D59").Value2 = "Fixation_D86_l.png"
$ws.Range("E59").Value2 = "Fixation_D86_r.png"
$ws.Range("H59").Value2 = "D86"
$ws.Range("B60").Value2 = "Face12_D86_S31"
$ws.Range("D60").Value2 = "Face12_D86_S31_l.png"
$ws.Range("E60").Value2 = "Face12_D86_S31_r.png"
$ws.Range("H60").Value2 = "D86"
$ws.Range("J60").Value2 = "S31"
$ws.Range("D61").Value2 = "Fixation_D86_l.png"
$ws.Range("E61").Value2 = "Fixation_D86_r.png"
$ws.Range("H61").Value2 = "D86"
$ws.Range("B62").Value2 = "Face08_D86_S25"
$ws.Range("D62").Value2 = "Face08_D86_S25_l.png"
$ws.Range("E62").Value2 = "Face08_D86_S25_r.png"
$ws.Range("H62").Value2 = "D86"
$ws.Range("D63").Value2 = "Fixation_D86_l.png"
$ws.Range("E63").Value2 = "Fixation_D86_r.png"
$ws.Range("H63").Value2 = "D86"
$ws.Range("B64").Value2 = "Face18_D86_S25"
$ws.Range("D64").Value2 = "Face18_D86_S25_l.png"
$ws.Range("E64").Value2 = "Face18_D86_S25_r.png"
$ws.Range("H64").Value2 = "D86"
$ws.Range("D65").Value2 = "Fixation_D86_l.png"
$ws.Range("E65").Value2 = "Fixation_D86_r.png"
$ws.Range("H65").Value2 = "D86"
$ws.Range("B66").Value2 = "Face15_D86_S31"
$ws.Range("D66").Value2 = "Face15_D86_S31_l.png"
$ws.Range("E66").Value2 = "Face15_D86_S31_r.png"
$ws.Range("H66").Value2 = "D86"
$ws.Range("J66").Value2 = "S31"
$ws.Range("D67").Value2 = "Fixation_D86_l.png"
$ws.Range("E67").Value2 = "Fixation_D86_r.png"
$ws.Range("H67").Value2 = "D86"
$ws.Range("D68").Value2 = "Fixation_D55_l.png"
$ws.Range("E68").Value2 = "Fixation_D55_r.png"
$ws.Range("H68").Value2 = "D55"
$ws.Range("B69").Value2 = "Face16_D55_S25"
$ws.Range("D69").Value2 = "Face16_D55_S25_l.png"
$ws.Range("E69").Value2 = "Face16_D55_S25_r.png"
$ws.Range("H69").Value2 = "D55"
$ws.Range("D70").Value2 = "Fixation_D55_l.png"
$ws.Range("E70").Value2 = "Fixation_D55_r.png"
$ws.Range("H70").Value2 = "D55"
$ws.Range("B71").Value2 = "Face09_D55_S31"
$ws.Range("D71").Value2 = "Face09_D55_S31_l.png"
$ws.Range("E71").Value2 = "Face09_D55_S31_r.png"
$ws.Range("H71").Value2 = "D55"
$ws.Range("J71").Value2 = "S31"
$ws.Range("D72").Value2 = "Fixation_D55_l.png"
$ws.Range("E72").Value2 = "Fixation_D55_r.png"
$ws.Range("H72").Value2 = "D55"
$ws.Range("B73").Value2 = "Face06_D55_S20"
$ws.Range("D73").Value2 = "Face06_D55_S20_l.png"
$ws.Range("E73").Value2 = "Face06_D55_S20_r.png"
$ws.Range("H73").Value2 = "D55"
$ws.Range("D74").Value2 = "Fixation_D55_l.png"
$ws.Range("E74").Value2 = "Fixation_D55_r.png"
$ws.Range("H74").Value2 = "D55"
$ws.Range("B75").Value2 = "Face14_D55_S25"
$ws.Range("D75").Value2 = "Face14_D55_S25_l.png"
$ws.Range("E75").Value2 = "Face14_D55_S25_r.png"
$ws.Range("H75").Value2 = "D55"
$ws.Range("D76").Value2 = "Fixation_D55_l.png"
$ws.Range("E76").Value2 = "Fixation_D55_r.png"
$ws.Range("H76").Value2 = "D55"
$ws.Range("B77").Value2 = "Face01_D55_S20"
$ws.Range("D77").Value2 = "Face01_D55_S20_l.png"
$ws.Range("E77").Value2 = "Face01_D55_S20_r.png"
$ws.Range("H77").Value2 = "D55"
$ws.Range("D78").Value2 = "Fixation_D55_l.png"
$ws.Range("E78").Value2 = "Fixation_D55_r.png"
$ws.Range("H78").Value2 = "D55"
$ws.Range("B79").Value2 = "Face10_D55_S31"
$ws.Range("D79").Value2 = "Face10_D55_S31_l.png"
$ws.Range("E79").Value2 = "Face10_D55_S31_r.png"
$ws.Range("H79").Value2 = "D55"
$ws.Range("J79").Value2 = "S31"
$ws.Range("D80").Value2 = "Fixation_D55_l.png"
$ws.Range("E80").Value2 = "Fixation_D55_r.png"
$ws.Range("H80").Value2 = "D55"
$ws.Range("D81").Value2 = "Fixation_D69_l.png"
$ws.Range("E81").Value2 = "Fixation_D69_r.png"
$ws.Range("H81").Value2 = "D69"
$ws.Range("B82").Value2 = "Face13_D69_S20"
$ws.Range("D82").Value2 = "Face13_D69_S20_l.png"
$ws.Range("E82").Value2 = "Face13_D69_S20_r.png"
$ws.Range("H82").Value2 = "D69"
$ws.Range("D83").Value2 = "Fixation_D69_l.png"
$ws.Range("E83").Value2 = "Fixation_D69_r.png"
$ws.Range("H83").Value2 = "D69"
$ws.Range("B84").Value2 = "Face15_D69_S31"
$ws.Range("D84").Value2 = "Face15_D69_S31_l.png"
$ws.Range("E84").Value2 = "Face15_D69_S31_r.png"
$ws.Range("H84").Value2 = "D69"
$ws.Range("J84").Value2 = "S31"
$ws.Range("D85").Value2 = "Fixation_D69_l.png"
$ws.Range("E85").Value2 = "Fixation_D69_r.png"
$ws.Range("H85").Value2 = "D69"
$ws.Range("B86").Value2 = "Face04_D69_S25"
$ws.Range("D86").Value2 = "Face04_D69_S25_l.png"
$ws.Range("E86").Value2 = "Face04_D69_S25_r.png"
$ws.Range("H86").Value2 = "D69"
$ws.Range("D87").Value2 = "Fixation_D69_l.png"
$ws.Range("E87").Value2 = "Fixation_D69_r.png"
$ws.Range("H87").Value2 = "D69"
$ws.Range("B88").Value2 = "Face17_D69_S25"
$ws.Range("D88").Value2 = "Face17_D69_S25_l.png"
$ws.Range("E88").Value2 = "Face17_D69_S25_r.png"
$ws.Range("H88").Value2 = "D69"
$ws.Range("D89").Value2 = "Fixation_D69_l.png"
$ws.Range("E89").Value2 = "Fixation_D69_r.png"
$ws.Range("H89").Value2 = "D69"
$ws.Range("B90").Value2 = "Face06_D69_S20"
$ws.Range("D90").Value2 = "Face06_D69_S20_l.png"
$ws.Range("E90").Value2 = "Face06_D69_S20_r.png"
$ws.Range("H90").Value2 = "D69"
$ws.Range("D91").Value2 = "Fixation_D69_l.png"
$ws.Range("E91").Value2 = "Fixation_D69_r.png"
$ws.Range("H91").Value2 = "D69"
$ws.Range("B92").Value2 = "Face02_D69_S31"
$ws.Range("D92").Value2 = "Face02_D69_S31_l.png"
$ws.Range("E92").Value2 = "Face02_D69_S31_r.png"
$ws.Range("H92").Value2 = "D69"
$ws.Range("J92").Value2 = "S31"
$ws.Range("D93").Value2 = "Fixation_D69_l.png"
$ws.Range("E93").Value2 = "Fixation_D69_r.png"
$ws.Range("H93").Value2 = "D69"
$ws.Range("D94").Value2 = "Fixation_D55_l.png"
$ws.Range("E94").Value2 = "Fixation_D55_r.png"
$ws.Range("H94").Value2 = "D55"
$ws.Range("B95").Value2 = "Face03_D55_S25"
$ws.Range("D95").Value2 = "Face03_D55_S25_l.png"
$ws.Range("E95").Value2 = "Face03_D55_S25_r.png"
$ws.Range("H95").Value2 = "D55"
$ws.Range("D96").Value2 = "Fixation_D55_l.png"
$ws.Range("E96").Value2 = "Fixation_D55_r.png"
$ws.Range("H96").Value2 = "D55"
$ws.Range("B97").Value2 = "Face12_D55_S25"
$ws.Range("D97").Value2 = "Face12_D55_S25_l.png"
$ws.Range("E97").Value2 = "Face12_D55_S25_r.png"
$ws.Range("H97").Value2 = "D55"
$ws.Range("D98").Value2 = "Fixation_D55_l.png"
$ws.Range("E98").Value2 = "Fixation_D55_r.png"
$ws.Range("H98").Value2 = "D55"
$ws.Range("B99").Value2 = "Face18_D55_S31"
$ws.Range("D99").Value2 = "Face18_D55_S31_l.png"
$ws.Range("E99").Value2 = "Face18_D55_S31_r.png"
$ws.Range("H99").Value2 = "D55"
$ws.Range("J99").Value2 = "S31"
$ws.Range("D100").Value2 = "Fixation_D55_l.png"
$ws.Range("E100").Value2 = "Fixation_D55_r.png"
$ws.Range("H100").Value2 = "D55"
$ws.Range("B101").Value2 = "Face14_D55_S31"
$ws.Range("D101").Value2 = "Face14_D55_S31_l.png"
$ws.Range("E101").Value2 = "Face14_D55_S31_r.png"
$ws.Range("H101").Value2 = "D55"
$ws.Range("J101").Value2 = "S31"
$ws.Range("D102").Value2 = "Fixation_D55_l.png"
$ws.Range("E102").Value2 = "Fixation_D55_r.png"
$ws.Range("H102").Value2 = "D55"
$ws.Range("B103").Value2 = "Face05_D55_S20"
$ws.Range("D103").Value2 = "Face05_D55_S20_l.png"
$ws.Range("E103").Value2 = "Face05_D55_S20_r.png"
$ws.Range("H103").Value2 = "D55"
$ws.Range("D104").Value2 = "Fixation_D55_l.png"
$ws.Range("E104").Value2 = "Fixation_D55_r.png"
$ws.Range("H104").Value2 = "D55"
$ws.Range("B105").Value2 = "Face07_D55_S20"
$ws.Range("D105").Value2 = "Face07_D55_S20_l.png"
$ws.Range("E105").Value2 = "Face07_D55_S20_r.png"
$ws.Range("H105").Value2 = "D55"
$ws.Range("D106").Value2 = "Fixation_D55_l.png"
$ws.Range("E106").Value2 = "Fixation_D55_r.png"
$ws.Range("H106").Value2 = "D55"
$ws.Range("D107").Value2 = "Fixation_D86_l.png"
$ws.Range("E107").Value2 = "Fixation_D86_r.png"
$ws.Range("H107").Value2 = "D86"
$ws.Range("B108").Value2 = "Face16_D86_S25"
$ws.Range("D108").Value2 = "Face16_D86_S25_l.png"
$ws.Range("E108").Value2 = "Face16_D86_S25_r.png"
$ws.Range("H108").Value2 = "D86"
$ws.Range("D109").Value2 = "Fixation_D86_l.png"
$ws.Range("E109").Value2 = "Fixation_D86_r.png"
$ws.Range("H109").Value2 = "D86"
$ws.Range("B110").Value2 = "Face01_D86_S25"
$ws.Range("D110").Value2 = "Face01_D86_S25_l.png"
$ws.Range("E110").Value2 = "Face01_D86_S25_r.png"
$ws.Range("H110").Value2 = "D86"
$ws.Range("D111").Value2 = "Fixation_D86_l.png"
$ws.Range("E111").Value2 = "Fixation_D86_r.png"
$ws.Range("H111").Value2 = "D86"
$ws.Range("B112").Value2 = "Face11_D86_S20"
$ws.Range("D112").Value2 = "Face11_D86_S20_l.png"
$ws.Range("E112").Value2 = "Face11_D86_S20_r.png"
$ws.Range("H112").Value2 = "D86"
$ws.Range("D113").Value2 = "Fixation_D86_l.png"
$ws.Range("E113").Value2 = "Fixation_D86_r.png"
$ws.Range("H113").Value2 = "D86"
$ws.Range("B114").Value2 = "Face10_D86_S31"
$ws.Range("D114").Value2 = "Face10_D86_S31_l.png"
$ws.Range("E114").Value2 = "Face10_D86_S31_r.png"
$ws.Range("H114").Value2 = "D86"
$ws.Range("J114").Value2 = "S31"
$ws.Range("D115").Value2 = "Fixation_D86_l.png"
$ws.Range("E115").Value2 = "Fixation_D86_r.png"
$ws.Range("H115").Value2 = "D86"
$ws.Range("B116").Value2 = "Face09_D86_S31"
$ws.Range("D116").Value2 = "Face09_D86_S31_l.png"
$ws.Range("E116").Value2 = "Face09_D86_S31_r.png"
$ws.Range("H116").Value2 = "D86"
$ws.Range("J116").Value2 = "S31"
$ws.Range("D117").Value2 = "Fixation_D86_l.png"
$ws.Range("E117").Value2 = "Fixation_D86_r.png"
$ws.Range("H117").Value2 = "D86"
$ws.Range("B118").Value2 = "Face08_D86_S20"
$ws.Range("D118").Value2 = "Face08_D86_S20_l.png"
$ws.Range("E118").Value2 = "Face08_D86_S20_r.png"
$ws.Range("H118").Value2 = "D86"
$ws.Range("D119").Value2 = "Fixation_D86_l.png"
$ws.Range("E119").Value2 = "Fixation_D86_r.png"
$ws.Range("H119").Value2 = "D86"
$ws.Range("D120").Value2 = "Fixation_D69_l.png"
$ws.Range("E120").Value2 = "Fixation_D69_r.png"
$ws.Range("H120").Value2 = "D69"
$ws.Range("B121").Value2 = "Face06_D69_S31"
$ws.Range("D121").Value2 = "Face06_D69_S31_l.png"
$ws.Range("E121").Value2 = "Face06_D69_S31_r.png"
$ws.Range("H121").Value2 = "D69"
$ws.Range("J121").Value2 = "S31"
$ws.Range("D122").Value2 = "Fixation_D69_l.png"
$ws.Range("E122").Value2 = "Fixation_D69_r.png"
$ws.Range("H122").Value2 = "D69"
$ws.Range("B123").Value2 = "Face01_D69_S20"
$ws.Range("D123").Value2 = "Face01_D69_S20_l.png"
$ws.Range("E123").Value2 = "Face01_D69_S20_r.png"
$ws.Range("H123").Value2 = "D69"
$ws.Range("D124").Value2 = "Fixation_D69_l.png"
$ws.Range("E124").Value2 = "Fixation_D69_r.png"
$ws.Range("H124").Value2 = "D69"
$ws.Range("B125").Value2 = "Face05_D69_S25"
$ws.Range("D125").Value2 = "Face05_D69_S25_l.png"
$ws.Range("E125").Value2 = "Face05_D69_S25_r.png"
$ws.Range("H125").Value2 = "D69"
$ws.Range("D126").Value2 = "Fixation_D69_l.png"
$ws.Range("E126").Value2 = "Fixation_D69_r.png"
$ws.Range("H126").Value2 = "D69"
$ws.Range("B127").Value2 = "Face12_D69_S31"
$ws.Range("D127").Value2 = "Face12_D69_S31_l.png"
$ws.Range("E127").Value2 = "Face12_D69_S31_r.png"
$ws.Range("H127").Value2 = "D69"
$ws.Range("J127").Value2 = "S31"
$ws.Range("D128").Value2 = "Fixation_D69_l.png"
$ws.Range("E128").Value2 = "Fixation_D69_r.png"
$ws.Range("H128").Value2 = "D69"
$ws.Range("B129").Value2 = "Face09_D69_S20"
$ws.Range("D129").Value2 = "Face09_D69_S20_l.png"
$ws.Range("E129").Value2 = "Face09_D69_S20_r.png"
$ws.Range("H129").Value2 = "D69"
$ws.Range("D130").Value2 = "Fixation_D69_l.png"
$ws.Range("E130").Value2 = "Fixation_D69_r.png"
$ws.Range("H130").Value2 = "D69"
$ws.Range("B131").Value2 = "Face16_D69_S25"
$ws.Range("D131").Value2 = "Face16_D69_S25_l.png"
$ws.Range("E131").Value2 = "Face16_D69_S25_r.png"
$ws.Range("H131").Value2 = "D69"
$ws.Range("D132").Value2 = "Fixation_D69_l.png"
$ws.Range("E132").Value2 = "Fixation_D69_r.png"
$ws.Range("H132").Value2 = "D69"
$ws.Range("D133").Value2 = "Fixation_D86_l.png"
$ws.Range("E133").Value2 = "Fixation_D86_r.png"
$ws.Range("H133").Value2 = "D86"
$ws.Range("B134").Value2 = "Face15_D86_S25"
$ws.Range("D134").Value2 = "Face15_D86_S25_l.png"
$ws.Range("E134").Value2 = "Face15_D86_S25_r.png"
$ws.Range("H134").Value2 = "D86"
$ws.Range("D135").Value2 = "Fixation_D86_l.png"
$ws.Range("E135").Value2 = "Fixation_D86_r.png"
$ws.Range("H135").Value2 = "D86"
$ws.Range("B136").Value2 = "Face14_D86_S31"
$ws.Range("D136").Value2 = "Face14_D86_S31_l.png"
$ws.Range("E136").Value2 = "Face14_D86_S31_r.png"
$ws.Range("H136").Value2 = "D86"
$ws.Range("J136").Value2 = "S31"
$ws.Range("D137").Value2 = "Fixation_D86_l.png"
$ws.Range("E137").Value2 = "Fixation_D86_r.png"
$ws.Range("H137").Value2 = "D86"
$ws.Range("B138").Value2 = "Face11_D86_S31"
$ws.Range("D138").Value2 = "Face11_D86_S31_l.png"
$ws.Range("E138").Value2 = "Face11_D86_S31_r.png"
$ws.Range("H138").Value2 = "D86"
$ws.Range("J138").Value2 = "S31"
$ws.Range("D139").Value2 = "Fixation_D86_l.png"
$ws.Range("E139").Value2 = "Fixation_D86_r.png"
$ws.Range("H139").Value2 = "D86"
$ws.Range("B140").Value2 = "Face02_D86_S20"
$ws.Range("D140").Value2 = "Face02_D86_S20_l.png"
$ws.Range("E140").Value2 = "Face02_D86_S20_r.png"
$ws.Range("H140").Value2 = "D86"
$ws.Range("D141").Value2 = "Fixation_D86_l.png"
$ws.Range("E141").Value2 = "Fixation_D86_r.png"
$ws.Range("H141").Value2 = "D86"
$ws.Range("B142").Value2 = "Face18_D86_S20"
$ws.Range("D142").Value2 = "Face18_D86_S20_l.png"
$ws.Range("E142").Value2 = "Face18_D86_S20_r.png"
$ws.Range("H142").Value2 = "D86"
$ws.Range("D143").Value2 = "Fixation_D86_l.png"
$ws.Range("E143").Value2 = "Fixation_D86_r.png"
$ws.Range("H143").Value2 = "D86"
$ws.Range("B144").Value2 = "Face17_D86_S25"
$ws.Range("D144").Value2 = "Face17_D86_S25_l.png"
$ws.Range("E144").Value2 = "Face17_D86_S25_r.png"
$ws.Range("H144").Value2 = "D86"
$ws.Range("D145").Value2 = "Fixation_D86_l.png"
$ws.Range("E145").Value2 = "Fixation_D86_r.png"
$ws.Range("H145").Value2 = "D86"
$ws.Range("D146").Value2 = "Fixation_D55_l.png"
$ws.Range("E146").Value2 = "Fixation_D55_r.png"
$ws.Range("H146").Value2 = "D55"
$ws.Range("B147").Value2 = "Face03_D55_S31"
$ws.Range("D147").Value2 = "Face03_D55_S31_l.png"
$ws.Range("E147").Value2 = "Face03_D55_S31_r.png"
$ws.Range("H147").Value2 = "D55"
$ws.Range("J147").Value2 = "S31"
$ws.Range("D148").Value2 = "Fixation_D55_l.png"
$ws.Range("E148").Value2 = "Fixation_D55_r.png"
$ws.Range("H148").Value2 = "D55"
$ws.Range("B149").Value2 = "Face04_D55_S25"
$ws.Range("D149").Value2 = "Face04_D55_S25_l.png"
$ws.Range("E149").Value2 = "Face04_D55_S25_r.png"
$ws.Range("H149").Value2 = "D55"
$ws.Range("D150").Value2 = "Fixation_D55_l.png"
$ws.Range("E150").Value2 = "Fixation_D55_r.png"
$ws.Range("H150").Value2 = "D55"
$ws.Range("B151").Value2 = "Face10_D55_S20"
$ws.Range("D151").Value2 = "Face10_D55_S20_l.png"
$ws.Range("E151").Value2 = "Face10_D55_S20_r.png"
$ws.Range("H151").Value2 = "D55"
$ws.Range("D152").Value2 = "Fixation_D55_l.png"
$ws.Range("E152").Value2 = "Fixation_D55_r.png"
$ws.Range("H152").Value2 = "D55"
$ws.Range("B153").Value2 = "Face08_D55_S31"
$ws.Range("D153").Value2 = "Face08_D55_S31_l.png"
$ws.Range("E153").Value2 = "Face08_D55_S31_r.png"
$ws.Range("H153").Value2 = "D55"
$ws.Range("J153").Value2 = "S31"
$ws.Range("D154").Value2 = "Fixation_D55_l.png"
$ws.Range("E154").Value2 = "Fixation_D55_r.png"
$ws.Range("H154").Value2 = "D55"
$ws.Range("B155").Value2 = "Face07_D55_S25"
$ws.Range("D155").Value2 = "Face07_D55_S25_l.png"
$ws.Range("E155").Value2 = "Face07_D55_S25_r.png"
$ws.Range("H155").Value2 = "D55"
$ws.Range("D156").Value2 = "Fixation_D55_l.png"
$ws.Range("E156").Value2 = "Fixation_D55_r.png"
$ws.Range("H156").Value2 = "D55"
$ws.Range("B157").Value2 = "Face13_D55_S20"
$ws.Range("D157").Value2 = "Face13_D55_S20_l.png"
$ws.Range("E157").Value2 = "Face13_D55_S20_r.png"
$ws.Range("H157").Value2 = "D55"
$ws.Range("D158").Value2 = "Fixation_D55_l.png"
$ws.Range("E158").Value2 = "Fixation_D55_r.png"
$ws.Range("H158").Value2 = "D55"
$ws.Range("D159").Value2 = "Fixation_D69_l.png"
$ws.Range("E159").Value2 = "Fixation_D69_r.png"
$ws.Range("H159").Value2 = "D69"
